$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New time value in D13 (5 hours -> 0.208333... of a day), same style/number
# format as the other "Total time" cells (D6, D9) -> h:mm:ss
$ws.Range("D13").Value = 0.20833333333333334
$ws.Range("D13").NumberFormat = "h:mm:ss"

# New quick-summary note in E13
$ws.Range("E13").Value = "Set up NetBSD on VM, wrote a makefile. Makefile had to be specific, took too much time. (I will probably work over 40 hours this quarter, but I will document all time spent). "

# Update the active selection to E13, matching the saved workbook state
$ws.Range("E13").Select()
